# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col D) and "Correspond Handback DateTime"
# (col G) values for the 9f4ad892... entry (row 3) on the "zh-cn" sheet, and for the
# 9f4ad892... entry (row 3) on the "de-de" sheet, reflecting newly generated handback
# report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-02-22 09:22:37"
$wsZhCn.Range("G3").Value = "2016-02-22 09:23:23"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-02-22 09:22:49"
$wsDeDe.Range("G3").Value = "2016-02-22 09:23:45"
